$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.975.72"
$ws.Range("E2").Value = "  -2.64%  "
$ws.Range("D3").Value = "1.860.95"
$ws.Range("E3").Value = "  -2.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5096"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3738"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07113"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8876"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.54"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.12%  "
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").Value = "1.853.59"
$ws.Range("E13").Value = "  -2.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.291"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008350"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.85%  "
$ws.Range("E18").Value = "  -2.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "27.030.01"
$ws.Range("E20").Value = "  -2.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.054"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("D22").Value = "2.086.18"
$ws.Range("E22").Value = "  -3.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.463"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.58%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.836"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.080"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.680"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.652"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09046"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05121"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.059"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.154"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7289"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02043"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.053"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.482"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.071"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5330"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.598"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.310"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1470"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4617"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.05%  "
$ws.Range("E48").Value = "  -4.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.566"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.90%  "
